$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated net_utility (B) and net_pop (C) values for rows 2-30 (years 1991-2019),
# reflecting a re-run of the isoelastic welfare results analysis.
$updates = @(
    @{ Row = 2; B = -3947627781824.52; C = 101325389141.245 }
    @{ Row = 3; B = -3866653192551.22; C = 99553880368.6495 }
    @{ Row = 4; B = -3825634364722.95; C = 98695032200.904 }
    @{ Row = 5; B = -3816204213985.85; C = 98605769880.2815 }
    @{ Row = 6; B = -3837443484769.44; C = 99277657100.309 }
    @{ Row = 7; B = -3862990776704.82; C = 100038584275.633 }
    @{ Row = 8; B = -3875443214064.11; C = 100474592703.049 }
    @{ Row = 9; B = -3884734785625.92; C = 101016014567.082 }
    @{ Row = 10; B = -3939251559375.61; C = 102462434696.177 }
    @{ Row = 11; B = -4072341909898.26; C = 105813793852.587 }
    @{ Row = 12; B = -4149575891964.92; C = 107813784445.418 }
    @{ Row = 13; B = -4239855920111.04; C = 110051302238.458 }
    @{ Row = 14; B = -4313096651440.28; C = 112025583530.39 }
    @{ Row = 15; B = -4370199981593.58; C = 113643712208.792 }
    @{ Row = 16; B = -4491472902294.32; C = 116800379402.818 }
    @{ Row = 17; B = -4572653463280.02; C = 119057434316.519 }
    @{ Row = 18; B = -4574268440965.17; C = 119398412273.803 }
    @{ Row = 19; B = -4701212260339.85; C = 122622332027.74 }
    @{ Row = 20; B = -4765514789927.1; C = 124249537506.757 }
    @{ Row = 21; B = -4952833978390.29; C = 128952448529.837 }
    @{ Row = 22; B = -4991363016324.61; C = 130237865264.127 }
    @{ Row = 23; B = -5200621341410.96; C = 135442403225.982 }
    @{ Row = 24; B = -5370585693587.06; C = 139712938430.927 }
    @{ Row = 25; B = -5519205710419.74; C = 143485293037.281 }
    @{ Row = 26; B = -5678412781166.54; C = 147395449972.367 }
    @{ Row = 27; B = -5792831501880.72; C = 150306607400.455 }
    @{ Row = 28; B = -5918048872120.55; C = 153376526311.904 }
    @{ Row = 29; B = -6080677909939.26; C = 157661388182.435 }
    @{ Row = 30; B = -6158219641794.82; C = 159822323183.053 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

